$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped cryptocurrency price / 1h-volume-change data.
# Every touched cell in the sheet is stored as text (the source scraper
# writes inline strings, not numbers), so each assignment below uses a
# string literal. For values that look like plain numbers (e.g. "578.20")
# a leading apostrophe is used to force Excel to keep them as text instead
# of silently re-parsing them into a number and dropping the trailing zero.

$ws.Range("D2").Value = "62.480.49"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "3.017.96"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'578.20"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").Value = "'149.16"
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("D8").Value = "'0.525"
$ws.Range("E8").Value = "  -3.02%  "
$ws.Range("D9").Value = "3.015.01"
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("D10").Value = "'0.150"
$ws.Range("E10").Value = "  -4.14%  "
$ws.Range("D11").Value = "'5.69"
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("D12").Value = "'0.443"
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("D13").Value = "'0.0000230"
$ws.Range("E13").Value = "  -4.01%  "
$ws.Range("D14").Value = "'35.41"
$ws.Range("E14").Value = "  -5.02%  "
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "3.524.08"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "'7.02"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("D18").Value = "62.447.35"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "3.023.29"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").Value = "'471.57"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").Value = "'14.02"
$ws.Range("E21").Value = "  -3.75%  "
$ws.Range("D22").Value = "'0.695"
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("D23").Value = "'7.40"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("D24").Value = "'2.35"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").Value = "'80.87"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").Value = "'12.42"
$ws.Range("E26").Value = "  -3.23%  "
$ws.Range("D27").Value = "'10.48"
$ws.Range("E27").Value = "  +5.19%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'7.20"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").Value = "'2.62"
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("D32").Value = "'2.17"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").Value = "'27.20"
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").Value = "'0.109"
$ws.Range("E34").Value = "  -4.27%  "
$ws.Range("D35").Value = "'1.04"
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").Value = "0.0₃0795"
$ws.Range("E36").Value = "  -5.91%  "
$ws.Range("D37").Value = "'5.80"
$ws.Range("E37").Value = "  -3.86%  "
$ws.Range("D38").Value = "'2.16"
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").Value = "'3.02"
$ws.Range("E39").Value = "  -10.23%  "
$ws.Range("D40").Value = "'50.11"
$ws.Range("E40").Value = "  -1.00%  "
$ws.Range("D41").Value = "'8.99"
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("D42").Value = "'419.66"
$ws.Range("E42").Value = "  -5.10%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.113"
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.280"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("D45").Value = "2.801.83"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "'0.0356"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").Value = "'38.13"
$ws.Range("E47").Value = "  -4.16%  "
$ws.Range("D48").Value = "'127.29"
$ws.Range("E48").Value = "  -2.86%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'24.75"
$ws.Range("E50").Value = "  -3.26%  "
$ws.Range("D51").Value = "'0.108"
$ws.Range("E51").Value = "  -1.33%  "
